$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 7630
$ws.Range("C2").Value = 4

# Delete rows 3 and 4 entirely (they are removed from the sheet)
$ws.Rows("3:4").Delete()
